# Caderno de anotações: insere coluna "id" no início, numera as linhas de
# dados (1..9) e remove a última linha (questão 614 / Gestão da Qualidade).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove a última linha de dados (antiga linha 11) antes de reindexar colunas.
$ws.Rows.Item(11).Delete()

# Insere uma nova coluna A, empurrando data/questao_id/disciplina/assunto/anotacao
# para B:F.
$ws.Columns.Item(1).Insert()

# Cabeçalho da nova coluna.
$ws.Range("A1").Value = "id"

# Numera as 9 linhas de dados restantes (linhas 2..10).
for ($i = 2; $i -le 10; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 1
}
